$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Test_0_002_d_16" test case (previously row 2) is dropped from the
# sheet. Deleting the row shifts "Test_20_normal" and "Test_0_01_d_16" up
# into rows 2 and 3, shrinking the used range from A1:L4 to A1:L3.
$ws.Rows.Item(2).Delete()

# Leave the selection where it was left after making this edit.
$ws.Range("E14").Select()
